$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header label change: "Market cap" -> "MarketCap"
$ws.Range("B1").Value = "MarketCap"

# Match the final selection recorded in the sheet (B1 was the last touched cell)
[void]$ws.Range("B1").Select()
